$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new value would otherwise be auto-parsed as a number
# (these are Price column entries that must remain plain text, matching the source data)
$textFormatCells = @("D5", "D6", "D10", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D29", "D31", "D32", "D34", "D37", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "60.554.55"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.610.56"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "511.56"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "154.22"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").Value = "2.622.26"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "6.77"
$ws.Range("E10").Value = "  +4.28%  "
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "3.068.14"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "60.501.29"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "21.58"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "2.621.17"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "4.75"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").Value = "353.20"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "10.57"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "6.17"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "60.66"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").Value = "0.0₃0841"
$ws.Range("E28").Value = "  -3.61%  "
$ws.Range("D29").Value = "7.33"
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "19.39"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "151.37"
$ws.Range("E32").Value = "  -3.65%  "
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("D34").Value = "5.80"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("D37").Value = "0.888"
$ws.Range("E37").Value = "  +4.87%  "
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").Value = "36.29"
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").Value = "0.842"
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("D41").Value = "3.75"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").Value = "290.80"
$ws.Range("E42").Value = "  -6.72%  "
$ws.Range("D43").Value = "0.625"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "0.996"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "0.0553"
$ws.Range("E46").Value = "  -4.60%  "
$ws.Range("D47").Value = "19.75"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "4.91"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "1.997.77"
$ws.Range("E51").Value = "  -3.34%  "
